$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.720414177533826
$ws.Range("D2").Value = 0.2556391291964246
$ws.Range("E2").Value = 0.1906768770469398
$ws.Range("F2").Value = 5.130640175018812
$ws.Range("G2").Value = 0.002656722786199202
$ws.Range("I2").Value = 1.838950559567508
$ws.Range("J2").Value = 0.2146172273668157
$ws.Range("L2").Value = 1.467174508564312
$ws.Range("M2").Value = 0.7833411743702641
$ws.Range("B3").Value = 1.656119734008541
$ws.Range("D3").Value = 0.2286217290598529
$ws.Range("E3").Value = 0.1660768893720856
$ws.Range("F3").Value = 5.08843259718833
$ws.Range("G3").Value = 0.002666407169980712
$ws.Range("I3").Value = 1.856752350397834
$ws.Range("J3").Value = 0.1867347197371316
$ws.Range("L3").Value = 1.410391977242938
$ws.Range("M3").Value = 0.7535588197654022
$ws.Range("B4").Value = 1.617502453288012
$ws.Range("D4").Value = 0.2122036632750337
$ws.Range("E4").Value = 0.1510194226303057
$ws.Range("F4").Value = 5.066091950222443
$ws.Range("G4").Value = 0.002672656551162513
$ws.Range("I4").Value = 1.868695312051365
$ws.Range("J4").Value = 0.169606629061505
$ws.Range("L4").Value = 1.376446356163939
$ws.Range("M4").Value = 0.7357146997863424
$ws.Range("B5").Value = 1.601981403824141
$ws.Range("D5").Value = 0.2055542126355761
$ws.Range("E5").Value = 0.1448940402937637
$ws.Range("F5").Value = 5.05787896226758
$ws.Range("G5").Value = 0.002675279776920032
$ws.Range("I5").Value = 1.873815549216999
$ws.Range("J5").Value = 0.1626237854456463
$ws.Range("L5").Value = 1.362842012534827
$ws.Range("M5").Value = 0.7285536105984676
$ws.Range("B6").Value = 1.599417169398066
$ws.Range("D6").Value = 0.2044525008370215
$ws.Range("E6").Value = 0.1438775352558821
$ws.Range("F6").Value = 5.056568786476817
$ws.Range("G6").Value = 0.002675719994675727
$ws.Range("I6").Value = 1.874681028379626
$ws.Range("J6").Value = 0.1614640801816449
$ws.Range("L6").Value = 1.360596771963628
$ws.Range("M6").Value = 0.7273711733311501
$ws.Range("B7").Value = 1.61729225766166
$ws.Range("D7").Value = 0.2121138227426513
$ws.Range("E7").Value = 0.1509367719559336
$ws.Range("F7").Value = 5.065977589732626
$ws.Range("G7").Value = 0.00267269161863071
$ws.Range("I7").Value = 1.86876334104327
$ws.Range("J7").Value = 0.1695124694280992
$ws.Range("L7").Value = 1.376261959971373
$ws.Range("M7").Value = 0.7356176761848729
$ws.Range("B8").Value = 1.698066732295615
$ws.Range("D8").Value = 0.246286943371274
$ws.Range("E8").Value = 0.1821843247831936
$ws.Range("F8").Value = 5.115340168877538
$ws.Range("G8").Value = 0.002659999244843655
$ws.Range("I8").Value = 1.844877704620856
$ws.Range("J8").Value = 0.2050044269085873
$ws.Range("L8").Value = 1.447403687709141
$ws.Range("M8").Value = 0.7729799625777574
$ws.Range("B9").Value = 1.863317857506331
$ws.Range("D9").Value = 0.3147443468033941
$ws.Range("E9").Value = 0.2438902338305695
$ws.Range("F9").Value = 5.240880509278497
$ws.Range("G9").Value = 0.00263749975258265
$ws.Range("I9").Value = 1.806128663734825
$ws.Range("J9").Value = 0.2745874686203251
$ws.Range("L9").Value = 1.594322113465978
$ws.Range("M9").Value = 0.8497950409376074
$ws.Range("B10").Value = 1.988967177205495
$ws.Range("D10").Value = 0.3660503961091024
$ws.Range("E10").Value = 0.2895711167915067
$ws.Range("F10").Value = 5.351178513035023
$ws.Range("G10").Value = 0.002622405438287998
$ws.Range("I10").Value = 1.782670175959083
$ws.Range("J10").Value = 0.325772487748111
$ws.Range("L10").Value = 1.706959725180013
$ws.Range("M10").Value = 0.9084550999932191
$ws.Range("B11").Value = 2.047064447341143
$ws.Range("D11").Value = 0.3896393261319986
$ws.Range("E11").Value = 0.3104464408892795
$ws.Range("F11").Value = 5.40539950534523
$ws.Range("G11").Value = 0.002615845924408952
$ws.Range("I11").Value = 1.773103490561226
$ws.Range("J11").Value = 0.349087928938701
$ws.Range("L11").Value = 1.759262487474814
$ws.Range("M11").Value = 0.9356383698185056
$ws.Range("B12").Value = 2.069200420650191
$ws.Range("D12").Value = 0.3986100640087216
$ws.Range("E12").Value = 0.3183665221715302
$ws.Range("F12").Value = 5.426523238668494
$ws.Range("G12").Value = 0.002613405793262906
$ws.Range("I12").Value = 1.769641151853328
$ws.Range("J12").Value = 0.3579226480383397
$ws.Range("L12").Value = 1.779224317406204
$ws.Range("M12").Value = 0.9460047383996368
$ws.Range("B13").Value = 2.064426989430899
$ws.Range("D13").Value = 0.3966763234338373
$ws.Range("E13").Value = 0.3166600988098196
$ws.Range("F13").Value = 5.42194740543033
$ws.Range("G13").Value = 0.002613929375612408
$ws.Range("I13").Value = 1.770379672268184
$ws.Range("J13").Value = 0.3560196617370934
$ws.Range("L13").Value = 1.774918197149418
$ws.Range("M13").Value = 0.9437689082235607
$ws.Range("B14").Value = 2.048882860712297
$ws.Range("D14").Value = 0.3903765771088104
$ws.Range("E14").Value = 0.3110977209281316
$ws.Range("F14").Value = 5.407125458721566
$ws.Range("G14").Value = 0.00261564429705872
$ws.Range("I14").Value = 1.772815419369913
$ws.Range("J14").Value = 0.3498146460365774
$ws.Range("L14").Value = 1.760901617023364
$ws.Range("M14").Value = 0.9364897553946037
$ws.Range("B15").Value = 2.039379343142798
$ws.Range("D15").Value = 0.3865228325013845
$ws.Range("E15").Value = 0.307692606172921
$ws.Range("F15").Value = 5.398123898293875
$ws.Range("G15").Value = 0.002616700432191975
$ws.Range("I15").Value = 1.77432831312359
$ws.Range("J15").Value = 0.3460146689081114
$ws.Range("L15").Value = 1.752336451797078
$ws.Range("M15").Value = 0.9320405553417714
$ws.Range("B16").Value = 1.985189341355465
$ws.Range("D16").Value = 0.3645140156714319
$ws.Range("E16").Value = 0.288208893602615
$ws.Range("F16").Value = 5.347717275051679
$ws.Range("G16").Value = 0.002622840266987935
$ws.Range("I16").Value = 1.783317738596693
$ws.Range("J16").Value = 0.3242494960892941
$ws.Range("L16").Value = 1.703563268050232
$ws.Range("M16").Value = 0.9066887225358045
$ws.Range("B17").Value = 1.952186518965618
$ws.Range("D17").Value = 0.3510778103673431
$ws.Range("E17").Value = 0.2762815529096656
$ws.Range("F17").Value = 5.317837471460507
$ws.Range("G17").Value = 0.002626685245838953
$ws.Range("I17").Value = 1.789116460941457
$ws.Range("J17").Value = 0.3109060741729479
$ws.Range("L17").Value = 1.673916950369517
$ws.Range("M17").Value = 0.8912645507844701
$ws.Range("B18").Value = 1.933292476535428
$ws.Range("D18").Value = 0.3433730124024805
$ws.Range("E18").Value = 0.2694301042243978
$ws.Range("F18").Value = 5.301031626929671
$ws.Range("G18").Value = 0.002628925688631005
$ws.Range("I18").Value = 1.792555620699829
$ws.Range("J18").Value = 0.3032341461906185
$ws.Range("L18").Value = 1.656965074039192
$ws.Range("M18").Value = 0.8824398370779249
$ws.Range("B19").Value = 1.926910422122432
$ws.Range("D19").Value = 0.3407682436353809
$ws.Range("E19").Value = 0.2671117978243416
$ws.Range("F19").Value = 5.295406471438014
$ws.Range("G19").Value = 0.002629689239783388
$ws.Range("I19").Value = 1.793737858563723
$ws.Range("J19").Value = 0.3006370179234352
$ws.Range("L19").Value = 1.651242528431908
$ws.Range("M19").Value = 0.879459960440343
$ws.Range("B20").Value = 1.955690580611645
$ws.Range("D20").Value = 0.3525056844835603
$ws.Range("E20").Value = 0.2775503134753023
$ws.Range("F20").Value = 5.320978806361751
$ws.Range("G20").Value = 0.002626272951219644
$ws.Range("I20").Value = 1.788488414231644
$ws.Range("J20").Value = 0.3123262001780347
$ws.Range("L20").Value = 1.677062490673791
$ws.Range("M20").Value = 0.8929016232951454
$ws.Range("B21").Value = 2.053444855795021
$ws.Range("D21").Value = 0.3922259126017877
$ws.Range("E21").Value = 0.3127311063009586
$ws.Range("F21").Value = 5.411462894934346
$ws.Range("G21").Value = 0.002615139396597033
$ws.Range("I21").Value = 1.772095617233838
$ws.Range("J21").Value = 0.3516370460450275
$ws.Range("L21").Value = 1.765014370604433
$ws.Range("M21").Value = 0.9386258392716087
$ws.Range("B22").Value = 2.118124888788543
$ws.Range("D22").Value = 0.4184086587806064
$ws.Range("E22").Value = 0.3358121877663933
$ws.Range("F22").Value = 5.474051356361997
$ws.Range("G22").Value = 0.00260811820665057
$ws.Range("I22").Value = 1.762317379202429
$ws.Range("J22").Value = 0.3773625729341461
$ws.Range("L22").Value = 1.823406234252786
$ws.Range("M22").Value = 0.9689332239654362
$ws.Range("B23").Value = 2.083531187186281
$ws.Range("D23").Value = 0.4044132510923362
$ws.Range("E23").Value = 0.3234848160397945
$ws.Range("F23").Value = 5.440327629693627
$ws.Range("G23").Value = 0.002611842302293347
$ws.Range("I23").Value = 1.767450115844809
$ws.Range("J23").Value = 0.3636288947773494
$ws.Range("L23").Value = 1.792157092651223
$ws.Range("M23").Value = 0.9527184776652575
$ws.Range("B24").Value = 1.954106146091476
$ws.Range("D24").Value = 0.3518600809194652
$ws.Range("E24").Value = 0.2769766890558429
$ws.Range("F24").Value = 5.319557449472001
$ws.Range("G24").Value = 0.002626459256317635
$ws.Range("I24").Value = 1.788772025974403
$ws.Range("J24").Value = 0.3116841632222247
$ws.Range("L24").Value = 1.675640104987338
$ws.Range("M24").Value = 0.8921613694234054
$ws.Range("B25").Value = 1.817873305011631
$ws.Range("D25").Value = 0.296056934676642
$ws.Range("E25").Value = 0.22714263938893
$ws.Range("F25").Value = 5.203786272845974
$ws.Range("G25").Value = 0.002643332748179515
$ws.Range("I25").Value = 1.815737149650822
$ws.Range("J25").Value = 0.2557578240414955
$ws.Range("L25").Value = 1.553765317977394
$ws.Range("M25").Value = 0.828628536793758
